$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: fill in the new task details
$ws.Range("B12").Value = "Phân lọai các câu hỏi trong chương trình"
$ws.Range("C12").Value = "Hoang"
$ws.Range("D12").Value = "17/11"
$ws.Range("E12").Value = 0
$ws.Range("G12").Value = "on processing"

# Row height for row 12
$ws.Rows.Item(12).RowHeight = 33

# Update the active selection to K12
$ws.Range("K12").Select()
